$d = $word.ActiveDocument

$replacements = @(
    @("675÷9=75, 0", "497÷4=124, 1"),
    @("761÷9=84, 5", "839÷8=104, 7"),
    @("300÷6=50, 0", "246÷3=82, 0"),
    @("966÷6=161, 0", "663÷5=132, 3"),
    @("244÷2=122, 0", "332÷4=83, 0"),
    @("776÷6=129, 2", "820÷2=410, 0"),
    @("504÷4=126, 0", "597÷8=74, 5"),
    @("337÷7=48, 1", "983÷3=327, 2"),
    @("345÷6=57, 3", "443÷5=88, 3"),
    @("553÷9=61, 4", "409÷8=51, 1"),
    @("707÷6=117, 5", "255÷2=127, 1"),
    @("356÷9=39, 5", "851÷3=283, 2"),
    @("736÷6=122, 4", "223÷5=44, 3"),
    @("889÷2=444, 1", "209÷2=104, 1"),
    @("597÷6=99, 3", "455÷3=151, 2"),
    @("786÷8=98, 2", "301÷2=150, 1"),
    @("701÷2=350, 1", "517÷4=129, 1"),
    @("637÷6=106, 1", "741÷7=105, 6"),
    @("650÷4=162, 2", "387÷6=64, 3"),
    @("103÷8=12, 7", "319÷7=45, 4"),
    @("623÷4=155, 3", "326÷5=65, 1"),
    @("698÷4=174, 2", "134÷8=16, 6"),
    @("525÷3=175, 0", "305÷8=38, 1"),
    @("116÷5=23, 1", "193÷6=32, 1"),
    @("327÷2=163, 1", "481÷4=120, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"
